$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells. Numeric-looking price values in column D are
# written with a leading apostrophe (forcing text) and then the style is
# reset to Normal so Excel does not silently convert them to numbers
# (which would lose formatting such as trailing zeros or switch to
# scientific notation).

$ws.Range("D2").Value = "60.905.67"
$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("D3").Value = "3.362.87"
$ws.Range("E3").Value = "  -0.08%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'571.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.15%  "
$ws.Range("D6").Value = "'135.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.51%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.361.78"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "'0.474"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.31%  "
$ws.Range("D10").Value = "'7.59"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.36%  "
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").Value = "'0.388"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.99%  "
$ws.Range("D13").Value = "3.938.13"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("E14").Value = "  +1.67%  "
$ws.Range("D15").Value = "'0.0000172"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.46%  "
$ws.Range("D16").Value = "3.369.30"
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("D17").Value = "'25.18"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.75%  "
$ws.Range("D18").Value = "61.110.21"
$ws.Range("E18").Value = "  -1.88%  "
$ws.Range("D19").Value = "'13.90"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.67%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'5.80"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.51%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'9.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.74%  "
$ws.Range("D22").Value = "'372.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").Value = "'0.567"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.81%  "
$ws.Range("D24").Value = "3.498.58"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'70.56"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.91%  "
$ws.Range("E27").Value = "  +11.25%  "
$ws.Range("E28").Value = "  +21.02%  "
$ws.Range("D29").Value = "'7.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +11.59%  "
$ws.Range("E30").Value = "  +0.14%  "
$ws.Range("D31").Value = "'8.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.04%  "
$ws.Range("D32").Value = "'2.14"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.17%  "
$ws.Range("D33").Value = "'0.154"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.77%  "
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "3.394.43"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").Value = "'23.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.37%  "
$ws.Range("D37").Value = "'5.55"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.80%  "
$ws.Range("D38").Value = "'6.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.89%  "
$ws.Range("D40").Value = "'163.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("D41").Value = "'0.0786"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +5.24%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("B43").Value = "ONDO"
$ws.Range("C43").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D43").Value = "'1.21"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.64%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").Value = "'4.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.50%  "
$ws.Range("E45").Value = "  -0.21%  "
$ws.Range("D46").Value = "'41.29"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.09%  "
$ws.Range("E47").Value = "  +4.58%  "
$ws.Range("D48").Value = "'23.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.87%  "
$ws.Range("E49").Value = "  +5.64%  "
$ws.Range("D50").Value = "'23.07"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +15.22%  "
$ws.Range("E51").Value = "  +15.72%  "
